$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (A2:T13), keep header row 1 and its shared strings intact
$ws.Range("A2:T13").Delete()

# Rebuild shared-string order: column A (sending cluster) first for all rows,
# then column B, then C, then D -- matches how the authoring tool regenerated the sheet

# Column A
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(8,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,1).Value = "Resolving-Mac"

# Column B
$ws.Cells.Item(2,2).Value = "Podxl2"
$ws.Cells.Item(3,2).Value = "Podxl2"
$ws.Cells.Item(4,2).Value = "Podxl2"
$ws.Cells.Item(5,2).Value = "Podxl2"
$ws.Cells.Item(6,2).Value = "Podxl2"
$ws.Cells.Item(7,2).Value = "Podxl2"
$ws.Cells.Item(8,2).Value = "Podxl2"
$ws.Cells.Item(9,2).Value = "Podxl2"
$ws.Cells.Item(10,2).Value = "Podxl2"
$ws.Cells.Item(11,2).Value = "Podxl2"
$ws.Cells.Item(12,2).Value = "Podxl2"
$ws.Cells.Item(13,2).Value = "Podxl2"
$ws.Cells.Item(14,2).Value = "Podxl2"
$ws.Cells.Item(15,2).Value = "Podxl2"
$ws.Cells.Item(16,2).Value = "Podxl2"

# Column C
$ws.Cells.Item(2,3).Value = "Sell"
$ws.Cells.Item(3,3).Value = "Sell"
$ws.Cells.Item(4,3).Value = "Sell"
$ws.Cells.Item(5,3).Value = "Sell"
$ws.Cells.Item(6,3).Value = "Sell"
$ws.Cells.Item(7,3).Value = "Sell"
$ws.Cells.Item(8,3).Value = "Sell"
$ws.Cells.Item(9,3).Value = "Sell"
$ws.Cells.Item(10,3).Value = "Sell"
$ws.Cells.Item(11,3).Value = "Sell"
$ws.Cells.Item(12,3).Value = "Sell"
$ws.Cells.Item(13,3).Value = "Sell"
$ws.Cells.Item(14,3).Value = "Sell"
$ws.Cells.Item(15,3).Value = "Sell"
$ws.Cells.Item(16,3).Value = "Sell"

# Column D
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(3,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(6,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(12,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(15,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"

# Numeric columns E-T

# Column E
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(16,5).Value = 3

# Column F
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(16,6).Value = 1

# Column G
$ws.Cells.Item(2,7).Value = 1.860366666666667
$ws.Cells.Item(3,7).Value = 1.860366666666667
$ws.Cells.Item(4,7).Value = 1.860366666666667
$ws.Cells.Item(5,7).Value = 2.156435666666667
$ws.Cells.Item(6,7).Value = 2.156435666666667
$ws.Cells.Item(7,7).Value = 2.156435666666667
$ws.Cells.Item(8,7).Value = 0.1087846666666667
$ws.Cells.Item(9,7).Value = 0.1087846666666667
$ws.Cells.Item(10,7).Value = 0.1087846666666667
$ws.Cells.Item(11,7).Value = 0.446837
$ws.Cells.Item(12,7).Value = 0.446837
$ws.Cells.Item(13,7).Value = 0.446837
$ws.Cells.Item(14,7).Value = 0.1070173333333333
$ws.Cells.Item(15,7).Value = 0.1070173333333333
$ws.Cells.Item(16,7).Value = 0.1070173333333333

# Column H
$ws.Cells.Item(2,8).Value = 5.581099999999999
$ws.Cells.Item(3,8).Value = 5.581099999999999
$ws.Cells.Item(4,8).Value = 5.581099999999999
$ws.Cells.Item(5,8).Value = 6.469307000000001
$ws.Cells.Item(6,8).Value = 6.469307000000001
$ws.Cells.Item(7,8).Value = 6.469307000000001
$ws.Cells.Item(8,8).Value = 0.326354
$ws.Cells.Item(9,8).Value = 0.326354
$ws.Cells.Item(10,8).Value = 0.326354
$ws.Cells.Item(11,8).Value = 1.340511
$ws.Cells.Item(12,8).Value = 1.340511
$ws.Cells.Item(13,8).Value = 1.340511
$ws.Cells.Item(14,8).Value = 0.321052
$ws.Cells.Item(15,8).Value = 0.321052
$ws.Cells.Item(16,8).Value = 0.321052

# Column I
$ws.Cells.Item(2,9).Value = 0.3975617032346596
$ws.Cells.Item(3,9).Value = 0.3975617032346596
$ws.Cells.Item(4,9).Value = 0.3975617032346596
$ws.Cells.Item(5,9).Value = 0.4608318628349082
$ws.Cells.Item(6,9).Value = 0.4608318628349082
$ws.Cells.Item(7,9).Value = 0.4608318628349082
$ws.Cells.Item(8,9).Value = 0.0232473620070316
$ws.Cells.Item(9,9).Value = 0.0232473620070316
$ws.Cells.Item(10,9).Value = 0.0232473620070316
$ws.Cells.Item(11,9).Value = 0.09548939032893099
$ws.Cells.Item(12,9).Value = 0.09548939032893099
$ws.Cells.Item(13,9).Value = 0.09548939032893099
$ws.Cells.Item(14,9).Value = 0.02286968159446954
$ws.Cells.Item(15,9).Value = 0.02286968159446954
$ws.Cells.Item(16,9).Value = 0.02286968159446954

# Column J
$ws.Cells.Item(2,10).Value = 0.3975617032346596
$ws.Cells.Item(3,10).Value = 0.3975617032346596
$ws.Cells.Item(4,10).Value = 0.3975617032346596
$ws.Cells.Item(5,10).Value = 0.4608318628349083
$ws.Cells.Item(6,10).Value = 0.4608318628349083
$ws.Cells.Item(7,10).Value = 0.4608318628349083
$ws.Cells.Item(8,10).Value = 0.02324736200703161
$ws.Cells.Item(9,10).Value = 0.02324736200703161
$ws.Cells.Item(10,10).Value = 0.02324736200703161
$ws.Cells.Item(11,10).Value = 0.09548939032893101
$ws.Cells.Item(12,10).Value = 0.09548939032893101
$ws.Cells.Item(13,10).Value = 0.09548939032893101
$ws.Cells.Item(14,10).Value = 0.02286968159446954
$ws.Cells.Item(15,10).Value = 0.02286968159446954
$ws.Cells.Item(16,10).Value = 0.02286968159446954

# Column K
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(16,11).Value = 3

# Column L
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(16,12).Value = 1

# Column M
$ws.Cells.Item(2,13).Value = 0.353079
$ws.Cells.Item(3,13).Value = 24.359095
$ws.Cells.Item(4,13).Value = 0.6745613333333332
$ws.Cells.Item(5,13).Value = 0.353079
$ws.Cells.Item(6,13).Value = 24.359095
$ws.Cells.Item(7,13).Value = 0.6745613333333332
$ws.Cells.Item(8,13).Value = 0.353079
$ws.Cells.Item(9,13).Value = 24.359095
$ws.Cells.Item(10,13).Value = 0.6745613333333332
$ws.Cells.Item(11,13).Value = 0.353079
$ws.Cells.Item(12,13).Value = 24.359095
$ws.Cells.Item(13,13).Value = 0.6745613333333332
$ws.Cells.Item(14,13).Value = 0.353079
$ws.Cells.Item(15,13).Value = 24.359095
$ws.Cells.Item(16,13).Value = 0.6745613333333332

# Column N
$ws.Cells.Item(2,14).Value = 1.059237
$ws.Cells.Item(3,14).Value = 73.077285
$ws.Cells.Item(4,14).Value = 2.023684
$ws.Cells.Item(5,14).Value = 1.059237
$ws.Cells.Item(6,14).Value = 73.077285
$ws.Cells.Item(7,14).Value = 2.023684
$ws.Cells.Item(8,14).Value = 1.059237
$ws.Cells.Item(9,14).Value = 73.077285
$ws.Cells.Item(10,14).Value = 2.023684
$ws.Cells.Item(11,14).Value = 1.059237
$ws.Cells.Item(12,14).Value = 73.077285
$ws.Cells.Item(13,14).Value = 2.023684
$ws.Cells.Item(14,14).Value = 1.059237
$ws.Cells.Item(15,14).Value = 73.077285
$ws.Cells.Item(16,14).Value = 2.023684

# Column O
$ws.Cells.Item(2,15).Value = 0.01390801122570493
$ws.Cells.Item(3,15).Value = 0.9595205795530543
$ws.Cells.Item(4,15).Value = 0.02657140922124081
$ws.Cells.Item(5,15).Value = 0.01390801122570493
$ws.Cells.Item(6,15).Value = 0.9595205795530543
$ws.Cells.Item(7,15).Value = 0.02657140922124081
$ws.Cells.Item(8,15).Value = 0.01390801122570493
$ws.Cells.Item(9,15).Value = 0.9595205795530543
$ws.Cells.Item(10,15).Value = 0.02657140922124081
$ws.Cells.Item(11,15).Value = 0.01390801122570493
$ws.Cells.Item(12,15).Value = 0.9595205795530543
$ws.Cells.Item(13,15).Value = 0.02657140922124081
$ws.Cells.Item(14,15).Value = 0.01390801122570493
$ws.Cells.Item(15,15).Value = 0.9595205795530543
$ws.Cells.Item(16,15).Value = 0.02657140922124081

# Column P
$ws.Cells.Item(2,16).Value = 0.01390801122570493
$ws.Cells.Item(3,16).Value = 0.9595205795530543
$ws.Cells.Item(4,16).Value = 0.02657140922124081
$ws.Cells.Item(5,16).Value = 0.01390801122570493
$ws.Cells.Item(6,16).Value = 0.9595205795530543
$ws.Cells.Item(7,16).Value = 0.02657140922124081
$ws.Cells.Item(8,16).Value = 0.01390801122570493
$ws.Cells.Item(9,16).Value = 0.9595205795530543
$ws.Cells.Item(10,16).Value = 0.02657140922124081
$ws.Cells.Item(11,16).Value = 0.01390801122570493
$ws.Cells.Item(12,16).Value = 0.9595205795530543
$ws.Cells.Item(13,16).Value = 0.02657140922124081
$ws.Cells.Item(14,16).Value = 0.01390801122570493
$ws.Cells.Item(15,16).Value = 0.9595205795530543
$ws.Cells.Item(16,16).Value = 0.02657140922124081

# Column Q
$ws.Cells.Item(2,17).Value = 0.6568564022999999
$ws.Cells.Item(3,17).Value = 45.31684836816666
$ws.Cells.Item(4,17).Value = 1.254931419155555
$ws.Cells.Item(5,17).Value = 0.761392148751
$ws.Cells.Item(6,17).Value = 52.52882126572167
$ws.Cells.Item(7,17).Value = 1.454648118554222
$ws.Cells.Item(8,17).Value = 0.03840958132199999
$ws.Cells.Item(9,17).Value = 2.649896029876667
$ws.Cells.Item(10,17).Value = 0.07338192979288886
$ws.Cells.Item(11,17).Value = 0.157768761123
$ws.Cells.Item(12,17).Value = 10.884544932515
$ws.Cells.Item(13,17).Value = 0.3014189625026666
$ws.Cells.Item(14,17).Value = 0.037785573036
$ws.Cells.Item(15,17).Value = 2.606845389313333
$ws.Cells.Item(16,17).Value = 0.07218975506311111

# Column R
$ws.Cells.Item(2,18).Value = 5.911707620699999
$ws.Cells.Item(3,18).Value = 407.8516353135
$ws.Cells.Item(4,18).Value = 11.2943827724
$ws.Cells.Item(5,18).Value = 6.852529338759001
$ws.Cells.Item(6,18).Value = 472.7593913914951
$ws.Cells.Item(7,18).Value = 13.091833066988
$ws.Cells.Item(8,18).Value = 0.345686231898
$ws.Cells.Item(9,18).Value = 23.84906426889
$ws.Cells.Item(10,18).Value = 0.6604373681359998
$ws.Cells.Item(11,18).Value = 1.419918850107
$ws.Cells.Item(12,18).Value = 97.960904392635
$ws.Cells.Item(13,18).Value = 2.712770662524
$ws.Cells.Item(14,18).Value = 0.340070157324
$ws.Cells.Item(15,18).Value = 23.46160850382
$ws.Cells.Item(16,18).Value = 0.649707795568

# Column S
$ws.Cells.Item(2,19).Value = 0.005529292631498016
$ws.Cells.Item(3,19).Value = 0.3814686358958199
$ws.Cells.Item(4,19).Value = 0.01056377470734164
$ws.Cells.Item(5,19).Value = 0.006409254721470417
$ws.Cells.Item(6,19).Value = 0.4421776561038647
$ws.Cells.Item(7,19).Value = 0.01224495200957306
$ws.Cells.Item(8,19).Value = 0.0003233245717618217
$ws.Cells.Item(9,19).Value = 0.02230632226606662
$ws.Cells.Item(10,19).Value = 0.0006177151692031628
$ws.Cells.Item(11,19).Value = 0.001328067512630492
$ws.Cells.Item(12,19).Value = 0.09162403514958367
$ws.Cells.Item(13,19).Value = 0.00253728766671682
$ws.Cells.Item(14,19).Value = 0.0003180717883441797
$ws.Cells.Item(15,19).Value = 0.02194393013771923
$ws.Cells.Item(16,19).Value = 0.0006076796684061291

# Column T
$ws.Cells.Item(2,20).Value = 0.005529292631498018
$ws.Cells.Item(3,20).Value = 0.38146863589582
$ws.Cells.Item(4,20).Value = 0.01056377470734164
$ws.Cells.Item(5,20).Value = 0.006409254721470419
$ws.Cells.Item(6,20).Value = 0.4421776561038648
$ws.Cells.Item(7,20).Value = 0.01224495200957306
$ws.Cells.Item(8,20).Value = 0.0003233245717618219
$ws.Cells.Item(9,20).Value = 0.02230632226606662
$ws.Cells.Item(10,20).Value = 0.000617715169203163
$ws.Cells.Item(11,20).Value = 0.001328067512630492
$ws.Cells.Item(12,20).Value = 0.0916240351495837
$ws.Cells.Item(13,20).Value = 0.002537287666716821
$ws.Cells.Item(14,20).Value = 0.0003180717883441798
$ws.Cells.Item(15,20).Value = 0.02194393013771923
$ws.Cells.Item(16,20).Value = 0.0006076796684061291

Write-Output "done"